$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value2 = 256
$ws.Range("F3").Value2 = 1060
$ws.Range("F4").Value2 = 9418
$ws.Range("F5").Value2 = 198
$ws.Range("F6").Value2 = 70
$ws.Range("F7").Value2 = 6464
$ws.Range("F8").Value2 = 624
$ws.Range("F9").Value2 = 73
$ws.Range("F10").Value2 = 9884
$ws.Range("F11").Value2 = 11304
$ws.Range("F13").Value2 = 1166
$ws.Range("F14").Value2 = 4959
$ws.Range("F15").Value2 = 805
$ws.Range("F16").Value2 = 465
$ws.Range("F19").Value2 = 180
$ws.Range("F20").Value2 = 1343
$ws.Range("F21").Value2 = 255
$ws.Range("F22").Value2 = 1867
$ws.Range("F23").Value2 = 893
$ws.Range("F24").Value2 = 1262
$ws.Range("F25").Value2 = 858
$ws.Range("F27").Value2 = 2051
$ws.Range("F28").Value2 = 433
$ws.Range("F29").Value2 = 633
$ws.Range("F30").Value2 = 2690
$ws.Range("F31").Value2 = 187
$ws.Range("F32").Value2 = 1775
$ws.Range("F34").Value2 = 799
$ws.Range("F35").Value2 = 64
$ws.Range("F36").Value2 = 922
$ws.Range("F37").Value2 = 587
$ws.Range("F38").Value2 = 30
$ws.Range("F39").Value2 = 3352
$ws.Range("F40").Value2 = 238
$ws.Range("F41").Value2 = 86
$ws.Range("F42").Value2 = 518
$ws.Range("F43").Value2 = 582
$ws.Range("F44").Value2 = 28
$ws.Range("F45").Value2 = 899
$ws.Range("F46").Value2 = 241
$ws.Range("F47").Value2 = 6
$ws.Range("F48").Value2 = 4215
$ws.Range("F49").Value2 = 54

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value2 = 12
$ws.Range("F6").Value2 = 11
$ws.Range("F9").Value2 = 8
$ws.Range("F23").Value2 = 70

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value2 = 5966

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value2 = 256
$ws.Range("F3").Value2 = 1060
$ws.Range("F4").Value2 = 9418
$ws.Range("F5").Value2 = 198
$ws.Range("F6").Value2 = 70
$ws.Range("F7").Value2 = 12
$ws.Range("F8").Value2 = 6464
$ws.Range("F9").Value2 = 624
$ws.Range("F10").Value2 = 9884
$ws.Range("F11").Value2 = 11304
$ws.Range("F13").Value2 = 1166
$ws.Range("F14").Value2 = 4959
$ws.Range("F15").Value2 = 805
$ws.Range("F16").Value2 = 465
$ws.Range("F20").Value2 = 180
$ws.Range("F21").Value2 = 1343
$ws.Range("F22").Value2 = 255
$ws.Range("F23").Value2 = 1867
$ws.Range("F24").Value2 = 858
$ws.Range("F26").Value2 = 2051
$ws.Range("F27").Value2 = 433
$ws.Range("F28").Value2 = 633
$ws.Range("F29").Value2 = 2690
$ws.Range("F30").Value2 = 187
$ws.Range("F31").Value2 = 1776
$ws.Range("F34").Value2 = 799
$ws.Range("F39").Value2 = 64
$ws.Range("F40").Value2 = 922
$ws.Range("F41").Value2 = 587
$ws.Range("F42").Value2 = 31
$ws.Range("F44").Value2 = 238
$ws.Range("F45").Value2 = 582
$ws.Range("F46").Value2 = 899
$ws.Range("F47").Value2 = 241
$ws.Range("F48").Value2 = 6
$ws.Range("F49").Value2 = 4215
